# "First step to rebuild": fill in B1 with the same date as A1, seed row 2
# with a 0/100 pair (A2 -> 0, B2 -> 100, the 100 that used to live in A2),
# and update the selection to span A1:L2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 gets the same (numeric/date-serial) value currently stored in A1.
$ws.Range("B1").Value2 = $ws.Range("A1").Value2

# Row 2: the former A2 value (100) moves to B2, A2 becomes 0.
$ws.Range("B2").Value2 = 100
$ws.Range("A2").Value2 = 0

# Update the active selection to A1:L2.
$ws.Range("A1:L2").Select()
